# FINFLUX-3612 Cartias specific scenarios
#
# Inserts a new "Branchmanager" sheet (copy of the Login sheet's
# username/password layout, but with "Branchmanager" in B1 instead of the
# secret) right after "Login" and before "Verify", then refreshes the
# Login sheet's column width / selected cell.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")

# Insert the new worksheet immediately after "Login".
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "Branchmanager"

# Pick up the same header/value cell formatting (grey/green fills) that
# the Login sheet already uses, rather than re-creating new styles.
$loginSheet.Range("A1:B2").Copy()
$newSheet.Range("A1:B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "Branchmanager"
$newSheet.Range("A2").Value = "password"
$newSheet.Range("B2").Value = "password"

# Match column widths to their content.
$newSheet.Columns("A:B").AutoFit()
$newSheet.Range("E9").Select()

# Login sheet gets a refreshed auto-fit column A and a new selected cell.
$loginSheet.Columns("A").AutoFit()
$loginSheet.Range("D6").Select()

# Keep "Login" as the active/visible tab.
$loginSheet.Activate()
